$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# TC3 (rows 26-32): step 3 now is "excluir comprovante" instead of "detalhar solicitação de diária"
$ws.Range("B32").Value = "Chefe Clica em excluir comprovante."
$ws.Range("D32").Value = "SYSTEM Exclui o comprovante."

# TC4 (rows 35-41): step 3 now is "visualizar comprovante" instead of "excluir comprovante"
$ws.Range("B41").Value = "Chefe Clica em visualizar comprovante."
$ws.Range("D41").Value = "SYSTEM Exibe modal com o comprovante."

# TC5 (rows 44-50): step 3 now is "detalhar solicitação de diária" instead of "visualizar comprovante"
$ws.Range("B50").Value = "Chefe Clica para detalhar a solicitação de diária."
$ws.Range("D50").Value = "SYSTEM Apresenta a tela de Detalhar Diárias"

# TC8 (rows 69-74): expected result swapped with TC9's
$ws.Range("D74").Value = "SYSTEM Identifica que a prestação de contas indicada pelo usuário não está em nenhum desses dois estados: a) NÃO REALIZADA e b) DEVOLVIDA; Permite não permite um novo envio ou alterações na prestação (exclusão de documentos)."

# TC9 (rows 77-82): expected result swapped with TC8's
$ws.Range("D82").Value = "SYSTEM Identifica que a solicitação indicada pelo usuário ainda não pode ter sua prestação de contas realizada; Exibe mensagem de erro (MSG212 - Prestação de contas ainda não pode ser realizada) para o usuário, impedindo que ele preste contas (anexa arquivos e etc)."
